$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "-"
$ws.Range("C3").Value = "ELT-2A-Acionamentos"
$ws.Range("D3").Value = "-"
$ws.Range("F3").Value = "MCT-1A-Circuitos Elétricos"
$ws.Range("C6").Value = "MCT-2A-Acionamentos"
$ws.Range("D6").Value = "-"
$ws.Range("F6").Value = "MEC-1A-Circuitos Elétricos"
$ws.Range("F7").Value = "MEC-1A-Circuitos Elétricos"
